{"js": "// Each three-digit x one-digit multiplication fact in the document is\n// replaced with a new fact, matched by its exact original text.\nconst replacements = [\n  [\"929\u00d79=8361\", \"741\u00d77=5187\"],\n  [\"497\u00d74=1988\", \"646\u00d75=3230\"],\n  [\"601\u00d77=4207\", \"112\u00d74=448\"],\n  [\"709\u00d78=5672\", \"980\u00d79=8820\"],\n  [\"753\u00d74=3012\", \"203\u00d72=406\"],\n  [\"901\u00d76=5406\", \"503\u00d78=4024\"],\n  [\"608\u00d75=3040\", \"646\u00d73=1938\"],\n  [\"994\u00d72=1988\", \"671\u00d75=3355\"],\n  [\"912\u00d73=2736\", \"558\u00d78=4464\"],\n  [\"148\u00d72=296\", \"883\u00d73=2649\"],\n  [\"208\u00d78=1664\", \"447\u00d74=1788\"],\n  [\"651\u00d74=2604\", \"255\u00d73=765\"],\n  [\"424\u00d76=2544\", \"152\u00d73=456\"],\n  [\"597\u00d76=3582\", \"332\u00d73=996\"],\n  [\"626\u00d79=5634\", \"876\u00d78=7008\"],\n  [\"457\u00d77=3199\", \"980\u00d72=1960\"],\n  [\"529\u00d72=1058\", \"212\u00d76=1272\"],\n  [\"397\u00d75=1985\", \"985\u00d75=4925\"],\n  [\"558\u00d74=2232\", \"459\u00d72=918\"],\n  [\"441\u00d77=3087\", \"774\u00d79=6966\"],\n  [\"874\u00d79=7866\", \"747\u00d79=6723\"],\n  [\"582\u00d76=3492\", \"640\u00d79=5760\"],\n  [\"962\u00d77=6734\", \"998\u00d75=4990\"],\n  [\"800\u00d73=2400\", \"515\u00d76=3090\"],\n  [\"525\u00d74=2100\", \"543\u00d78=4344\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each three-digit x one-digit multiplication fact in the document is\n# replaced with a new fact, matched by its exact original text.\n$replacements = @(\n    @{Old=\"929\u00d79=8361\"; New=\"741\u00d77=5187\"},\n    @{Old=\"497\u00d74=1988\"; New=\"646\u00d75=3230\"},\n    @{Old=\"601\u00d77=4207\"; New=\"112\u00d74=448\"},\n    @{Old=\"709\u00d78=5672\"; New=\"980\u00d79=8820\"},\n    @{Old=\"753\u00d74=3012\"; New=\"203\u00d72=406\"},\n    @{Old=\"901\u00d76=5406\"; New=\"503\u00d78=4024\"},\n    @{Old=\"608\u00d75=3040\"; New=\"646\u00d73=1938\"},\n    @{Old=\"994\u00d72=1988\"; New=\"671\u00d75=3355\"},\n    @{Old=\"912\u00d73=2736\"; New=\"558\u00d78=4464\"},\n    @{Old=\"148\u00d72=296\"; New=\"883\u00d73=2649\"},\n    @{Old=\"208\u00d78=1664\"; New=\"447\u00d74=1788\"},\n    @{Old=\"651\u00d74=2604\"; New=\"255\u00d73=765\"},\n    @{Old=\"424\u00d76=2544\"; New=\"152\u00d73=456\"},\n    @{Old=\"597\u00d76=3582\"; New=\"332\u00d73=996\"},\n    @{Old=\"626\u00d79=5634\"; New=\"876\u00d78=7008\"},\n    @{Old=\"457\u00d77=3199\"; New=\"980\u00d72=1960\"},\n    @{Old=\"529\u00d72=1058\"; New=\"212\u00d76=1272\"},\n    @{Old=\"397\u00d75=1985\"; New=\"985\u00d75=4925\"},\n    @{Old=\"558\u00d74=2232\"; New=\"459\u00d72=918\"},\n    @{Old=\"441\u00d77=3087\"; New=\"774\u00d79=6966\"},\n    @{Old=\"874\u00d79=7866\"; New=\"747\u00d79=6723\"},\n    @{Old=\"582\u00d76=3492\"; New=\"640\u00d79=5760\"},\n    @{Old=\"962\u00d77=6734\"; New=\"998\u00d75=4990\"},\n    @{Old=\"800\u00d73=2400\"; New=\"515\u00d76=3090\"},\n    @{Old=\"525\u00d74=2100\"; New=\"543\u00d78=4344\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
